$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = "EloPP V1, Scores Logistic, Regularization = calibrated = 0.22"
$ws.Range("C10").Value = 0.40623008569885799
$ws.Range("C10").Style = "Good"

# Row 11
$ws.Range("A11").Value = "TeamSeedPredictor"
$ws.Range("C11").Value = 0.408240672723793

# Row 12
$ws.Range("A12").Value = "EloPP V1, Scores Logistic, Regularization = calibrated = 0.22, RandomSeed=421"
$ws.Range("C12").Value = 0.406812534700287
$ws.Range("C12").Style = "Note"

$ws.Range("C12").Select()
